$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: selection becomes the whole A1:C4 range (active cell = top-left A1)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$ws1.Range("A1:C4").Select()

# ---------------------------------------------------------------------------
# Sheet2: add a selection on C1 (previously no selection element at all)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("C1").Select()

# ---------------------------------------------------------------------------
# Sheet5: selection moves from B2 (range A1:C4) to just C1
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Sheet5")
$ws5.Activate()
$ws5.Range("C1").Select()

# ---------------------------------------------------------------------------
# Sheet6: selection moves from A4 to F25
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Sheet6")
$ws6.Activate()
$ws6.Range("F25").Select()

# ---------------------------------------------------------------------------
# Sheet9: C1 text value changes to match B1 ("UMAR--MZ007--LKJ--11--A"),
# and selection moves from B1 (range A1:C4) to just C1
# ---------------------------------------------------------------------------
$ws9 = $wb.Worksheets.Item("Sheet9")
$ws9.Activate()
$ws9.Range("C1").Value = "UMAR--MZ007--LKJ--11--A"
$ws9.Range("C1").Select()

# ---------------------------------------------------------------------------
# Sheet10: selection moves from A2 (range A2:C5) to just C2
# ---------------------------------------------------------------------------
$ws10 = $wb.Worksheets.Item("Sheet10")
$ws10.Activate()
$ws10.Range("C2").Select()

# ---------------------------------------------------------------------------
# Sheet11: rebuilt as a 4-column table (A1:D4), mirroring Sheet1's data with
# an extra column D holding a new "UMAR--MZ002--12dfg--M" series.
# ---------------------------------------------------------------------------
$ws11 = $wb.Worksheets.Item("Sheet11")
$ws11.Activate()

$ws11.Cells.Clear()

$ws11.Range("A1").Value = "period "
$ws11.Range("B1").Value = "UMAR--MZ002--1--M"
$ws11.Range("C1").Value = "UMAR--MZ002--12--M"
$ws11.Range("D1").Value = "UMAR--MZ002--12dfg--M"

$ws11.Range("A2").Value = "2020M01"
$ws11.Range("B2").Value = 11
$ws11.Range("C2").Value = 21
$ws11.Range("D2").Value = 21

$ws11.Range("A3").Value = "2020M02"
$ws11.Range("B3").Value = 212
$ws11.Range("C3").Value = 221
$ws11.Range("D3").Value = 221

$ws11.Range("A4").Value = "2020M03"
$ws11.Range("B4").Value = 213
$ws11.Range("C4").Value = 112
$ws11.Range("D4").Value = 112

$ws11.Range("D1:D4").Select()

# Leave Sheet11 as the active sheet/tab, matching activeTab="10" in workbook.xml
$ws11.Activate()
